$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "Disturbance, Momentum, State Observer에 대해서"
$ws.Range("E28").Value = "https://ropiens.tistory.com/179"

$ws.Range("D36").Value = "Uncertainty Quantification in Deep Learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/350"

$ws.Range("D52").Value = "[공지] 데이터사이언스 문제해결 모임"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2420&utm_source=rss&utm_medium=rss&utm_campaign=%25ea%25b3%25b5%25ec%25a7%2580-%25eb%258d%25b0%25ec%259d%25b4%25ed%2584%25b0%25ec%2582%25ac%25ec%259d%25b4%25ec%2596%25b8%25ec%258a%25a4-%25eb%25ac%25b8%25ec%25a0%259c%25ed%2595%25b4%25ea%25b2%25b0-%25eb%25aa%25a8%25ec%259e%2584"
